# Regenerate merged AHB file:
#   1. Rename the diff-header columns: "<Base>_old" -> "<Base>_FV2310"
#      and "<Base>_new" -> "<Base>_FV2404" (the "diff" column is untouched).
#   2. Turn the header + data range into a real Excel Table ("Table1").
#   3. Freeze the header row (split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (1-based): A..J = "<Base>_old", K = "diff", L..U = "<Base>_new"
$baseNames   = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols     = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$newCols     = @(12, 13, 14, 15, 16, 17, 18, 19, 20, 21)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $oldCols[$i]).Value = $baseNames[$i] + "_FV2310"
    $ws.Cells.Item(1, $newCols[$i]).Value = $baseNames[$i] + "_FV2404"
}

# Wrap the used range A1:U78 (header row + 77 data rows) in an Excel Table.
$tableRange = $ws.Range("A1:U78")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze panes below the header row (row 1), with focus back on the top-left
# data cell (A2), matching the saved "frozen" pane state.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
